$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ATOL 1e-5 -> 1e-4, metodo RK2 -> DOPRI5(4)
$ws.Range("C2").Value = 0.0001
$ws.Range("E2").Value = "DOPRI5(4)"

# Row 3: ATOL 1e-5 -> 1e-4, metodo RK2 -> DOPRI5(4)
$ws.Range("C3").Value = 0.0001
$ws.Range("E3").Value = "DOPRI5(4)"

# Row 4: ATOL 1e-5 -> 1e-6
$ws.Range("C4").Value = 0.000001

# Row 5: ATOL 1e-5 -> 1e-6
$ws.Range("C5").Value = 0.000001

# Row 6: ATOL 1e-7 -> 1e-6, RTOL 1e-7 -> 1e-6, metodo DOPRI5(4) -> SSPRK3
$ws.Range("C6").Value = 0.000001
$ws.Range("D6").Value = 0.000001
$ws.Range("E6").Value = "SSPRK3"

# Row 7: ATOL 1e-5 -> 1e-4, metodo RK2 -> DOPRI5(4)
$ws.Range("C7").Value = 0.0001
$ws.Range("E7").Value = "DOPRI5(4)"

# Row 8: ATOL 1e-7 -> 1e-6, RTOL 1e-7 -> 1e-6, metodo DOPRI5(4) -> SSPRK3
$ws.Range("C8").Value = 0.000001
$ws.Range("D8").Value = 0.000001
$ws.Range("E8").Value = "SSPRK3"

# Row 9: metodo ralston4 -> DOPRI5(4)
$ws.Range("E9").Value = "DOPRI5(4)"

# Row 10: ATOL 1e-5 -> 1e-4, metodo RK2 -> DOPRI5(4)
$ws.Range("C10").Value = 0.0001
$ws.Range("E10").Value = "DOPRI5(4)"

# Row 11: ATOL 1e-7 -> 1e-4, RTOL 1e-7 -> 1e-6
$ws.Range("C11").Value = 0.0001
$ws.Range("D11").Value = 0.000001

# Row 12: ATOL 1e-5 -> 1e-6
$ws.Range("C12").Value = 0.000001

# Row 13: metodo RK2 -> RK4
$ws.Range("E13").Value = "RK4"

# Row 14: ATOL 1e-5 -> 1e-6
$ws.Range("C14").Value = 0.000001

# Row 15: ATOL 1e-6 -> 1e-4, metodo SSPRK3 -> Fehlberg4(5)
$ws.Range("C15").Value = 0.0001
$ws.Range("E15").Value = "Fehlberg4(5)"

# Row 16: metodo RK2 -> Ralston3
$ws.Range("E16").Value = "Ralston3"

# Column E width widened to fit "Fehlberg4(5)"
$ws.Columns("E").ColumnWidth = 11.8

# Final selection left on C18 (below the table)
$ws.Range("C18").Select() | Out-Null
